$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update week 08 readings: the Friday reader for week 06 changes from
# "Samantha" to "Clara", and the Friday reader for week 08 changes from
# "Clara" to "Simmie".
$ws.Range("D6").Value = "Clara"
$ws.Range("D8").Value = "Simmie"

# Update the active cell selection to D14
$ws.Range("D14").Select()
